$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "(302962915, Asher  Odeh: -9,-6)"
$ws.Range("B1").Value = "(305251175, Or  Leder: 2,-6)"
$ws.Range("C1").Value = "(206532695, Matan Vakrat: 3,0)"
$ws.Range("D1").Value = "(203957296, Omri Ben Shabat: 1,4)"
$ws.Range("E1").Value = "(308035542, Anastasia  Kubi: 0,2)"
$ws.Range("F1").Value = "(308051846, Eyal  Sofer: -5,3)"
$ws.Range("G1").Value = "(311177802, Christina  Uksusman: -7,8)"

$ws.Range("A3").Value = "cost: 416.184492796085"
$ws.Range("A4").Value = "time: 48.898061599510626"
